# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders the worker/period rows in the account-statement table:
#   - Row 16 now holds EDITH TERESA URBINA MEZA's data (CC 45453825, period 1802),
#     with an updated Salario Basico (G16) of 1423500.
#   - Row 17 keeps ANDREA PAOLA GOMEZ ARRIETA's period 1811 entry unchanged.
#   - Row 18 now holds ANDREA PAOLA GOMEZ ARRIETA's period 1812 entry
#     (same Valor Mora / Salario Basico that used to sit in row 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: EDITH TERESA URBINA MEZA / period 1802 -----------------------
$ws.Range("C16").Value = "45453825"
$ws.Range("D16").Value = "EDITH TERESA URBINA MEZA"
$ws.Range("E16").Value = "1802"
$ws.Range("F16").Value = 16667
$ws.Range("G16").Value = 1423500

# --- Row 17: ANDREA PAOLA GOMEZ ARRIETA / period 1811 (unchanged) ---------
$ws.Range("C17").Value = "1143364030"
$ws.Range("D17").Value = "ANDREA PAOLA GOMEZ ARRIETA"
$ws.Range("E17").Value = "1811"
$ws.Range("F17").Value = 14933
$ws.Range("G17").Value = 1400000

# --- Row 18: ANDREA PAOLA GOMEZ ARRIETA / period 1812 ----------------------
$ws.Range("C18").Value = "1143364030"
$ws.Range("D18").Value = "ANDREA PAOLA GOMEZ ARRIETA"
$ws.Range("E18").Value = "1812"
$ws.Range("F18").Value = 56000
$ws.Range("G18").Value = 1400000
